$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4: function name, filename, and line number
$ws.Range("A4").Value = "FirmwareStatus"
$ws.Range("C4").Value = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/entservices-softwareupdate/FirmwareUpdate/FirmwareUpdateImplementation.cpp"
$ws.Range("D4").Value = 1436

# Move the active selection to A6 (matches saved view state)
$ws.Range("A6").Select()
